# Updates cryptos list price/volume/coin data (GitHub Actions scrape refresh).
# Commit: "Updated cryptos list on Thu Oct  5 15:19:42 UTC 2023 with GitHub Actions"
#
# All target cells hold plain text in the workbook (coin names, URLs, price
# strings such as "27.872.26", and padded percentage strings). Several price
# strings are valid numeric literals (e.g. "1.00", "0.526"), so each cell's
# number format is forced to Text ("@") before the value is written; this
# prevents Excel from auto-converting them into numbers and losing the
# original text formatting (trailing zeros, multi-dot grouping, padding).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.872.26' },
    @{ Cell = 'E2'; Value = '  +1.60%  ' },
    @{ Cell = 'D3'; Value = '1.633.27' },
    @{ Cell = 'E3'; Value = '  -0.20%  ' },
    @{ Cell = 'E4'; Value = '  +0.26%  ' },
    @{ Cell = 'D5'; Value = '212.60' },
    @{ Cell = 'E5'; Value = '  +0.01%  ' },
    @{ Cell = 'D6'; Value = '0.526' },
    @{ Cell = 'E6'; Value = '  -1.00%  ' },
    @{ Cell = 'E7'; Value = '  +0.21%  ' },
    @{ Cell = 'D8'; Value = '23.09' },
    @{ Cell = 'E8'; Value = '  +0.29%  ' },
    @{ Cell = 'D9'; Value = '0.262' },
    @{ Cell = 'E9'; Value = '  +2.14%  ' },
    @{ Cell = 'E10'; Value = '  +0.43%  ' },
    @{ Cell = 'D11'; Value = '0.0891' },
    @{ Cell = 'E11'; Value = '  +0.63%  ' },
    @{ Cell = 'D12'; Value = '1.864.16' },
    @{ Cell = 'E12'; Value = '  -0.18%  ' },
    @{ Cell = 'D13'; Value = '1.632.75' },
    @{ Cell = 'E13'; Value = '  -0.42%  ' },
    @{ Cell = 'D14'; Value = '4.05' },
    @{ Cell = 'E14'; Value = '  +0.80%  ' },
    @{ Cell = 'D15'; Value = '0.557' },
    @{ Cell = 'E15'; Value = '  -4.10%  ' },
    @{ Cell = 'D16'; Value = '64.59' },
    @{ Cell = 'E16'; Value = '  +0.69%  ' },
    @{ Cell = 'D17'; Value = '27.842.73' },
    @{ Cell = 'E17'; Value = '  +1.51%  ' },
    @{ Cell = 'D18'; Value = '231.44' },
    @{ Cell = 'E18'; Value = '  +0.99%  ' },
    @{ Cell = 'E19'; Value = '  +0.02%  ' },
    @{ Cell = 'D20'; Value = '7.61' },
    @{ Cell = 'E20'; Value = '  -0.07%  ' },
    @{ Cell = 'E21'; Value = '  +0.21%  ' },
    @{ Cell = 'D22'; Value = '4.30' },
    @{ Cell = 'E22'; Value = '  -0.05%  ' },
    @{ Cell = 'D23'; Value = '9.98' },
    @{ Cell = 'E23'; Value = '  +2.37%  ' },
    @{ Cell = 'D24'; Value = '2.09' },
    @{ Cell = 'E24'; Value = '  +6.57%  ' },
    @{ Cell = 'D25'; Value = '150.06' },
    @{ Cell = 'E25'; Value = '  +0.39%  ' },
    @{ Cell = 'D26'; Value = '6.92' },
    @{ Cell = 'E26'; Value = '  -1.19%  ' },
    @{ Cell = 'B28'; Value = 'BinanceUSD' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Cell = 'D28'; Value = '1.00' },
    @{ Cell = 'E28'; Value = '  +0.28%  ' },
    @{ Cell = 'B29'; Value = 'EthereumClassic' },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Cell = 'D29'; Value = '15.63' },
    @{ Cell = 'E29'; Value = '  +0.33%  ' },
    @{ Cell = 'D30'; Value = '1.18' },
    @{ Cell = 'E30'; Value = '  -0.27%  ' },
    @{ Cell = 'D31'; Value = '0.0483' },
    @{ Cell = 'E31'; Value = '  -0.98%  ' },
    @{ Cell = 'E32'; Value = '  +0.61%  ' },
    @{ Cell = 'D33'; Value = '1.470.72' },
    @{ Cell = 'E33'; Value = '  +3.54%  ' },
    @{ Cell = 'D34'; Value = '3.09' },
    @{ Cell = 'E34'; Value = '  -2.56%  ' },
    @{ Cell = 'E35'; Value = '  -2.40%  ' },
    @{ Cell = 'E36'; Value = '  +0.57%  ' },
    @{ Cell = 'D37'; Value = '0.566' },
    @{ Cell = 'E37'; Value = '  -0.98%  ' },
    @{ Cell = 'B38'; Value = 'VeChain' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D38'; Value = '0.0168' },
    @{ Cell = 'E38'; Value = '  +0.53%  ' },
    @{ Cell = 'B39'; Value = 'ARBITRUM' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Cell = 'D39'; Value = '0.876' },
    @{ Cell = 'E39'; Value = '  +0.02%  ' },
    @{ Cell = 'D40'; Value = '0.923' },
    @{ Cell = 'E40'; Value = '  +7.08%  ' },
    @{ Cell = 'D41'; Value = '69.36' },
    @{ Cell = 'E41'; Value = '  +6.97%  ' },
    @{ Cell = 'B42'; Value = 'WEMIXToken' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D42'; Value = '1.02' },
    @{ Cell = 'E42'; Value = '  -1.09%  ' },
    @{ Cell = 'B43'; Value = 'PaxDollar' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Cell = 'D43'; Value = '1.00' },
    @{ Cell = 'E43'; Value = '  +0.24%  ' },
    @{ Cell = 'D44'; Value = '2.46' },
    @{ Cell = 'E44'; Value = '  -1.41%  ' },
    @{ Cell = 'E45'; Value = '  +0.05%  ' },
    @{ Cell = 'D46'; Value = '5.39' },
    @{ Cell = 'E46'; Value = '  -2.43%  ' },
    @{ Cell = 'D47'; Value = '1.774.47' },
    @{ Cell = 'E47'; Value = '  -0.15%  ' },
    @{ Cell = 'D48'; Value = '1.70' },
    @{ Cell = 'E48'; Value = '  +2.50%  ' },
    @{ Cell = 'D49'; Value = '85.93' },
    @{ Cell = 'E49'; Value = '  +0.40%  ' },
    @{ Cell = 'B50'; Value = 'Algorand' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D50'; Value = '0.0989' },
    @{ Cell = 'E50'; Value = '  -0.17%  ' },
    @{ Cell = 'B51'; Value = 'EnergySwap' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D51'; Value = '7.79' },
    @{ Cell = 'E51'; Value = '  +1.10%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = '@'
    $cell.Value = $u.Value
}
